$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 14601
$ws.Range("I32").Value = 5800
$ws.Range("J32").Value = 19001.5
$ws.Range("K32").Value = 5800
$ws.Range("L32").Value = 19001.5
$ws.Range("M32").Value = -5474
$ws.Range("N32").Value = -19653.5

$ws.Range("H88").Value = 13077.333
$ws.Range("I88").Value = 3331.3333
$ws.Range("J88").Value = 16326
$ws.Range("K88").Value = 3331.3333
$ws.Range("L88").Value = 16326
$ws.Range("M88").Value = -2925.3333
$ws.Range("N88").Value = -17138

$ws.Range("H91").Value = 13077.333
$ws.Range("I91").Value = 3331.3333
$ws.Range("J91").Value = 16326
$ws.Range("K91").Value = 3331.3333
$ws.Range("L91").Value = 16326
$ws.Range("M91").Value = -1927.3333
$ws.Range("N91").Value = -19134

$ws.Range("H92").Value = 2233
$ws.Range("I92").Value = 505
$ws.Range("J92").Value = 9145
$ws.Range("K92").Value = 505
$ws.Range("L92").Value = 9145
$ws.Range("M92").Value = 743
$ws.Range("N92").Value = -11641

$ws.Range("H125").Value = 971.4286
$ws.Range("I125").Value = 1006.2
$ws.Range("K125").Value = 9055.800000000001
$ws.Range("M125").Value = -6595.800000000001

$ws.Range("H131").Value = 18523.75
$ws.Range("I131").Value = 16698.334
$ws.Range("K131").Value = 50095.00199999999
$ws.Range("M131").Value = -45055.00199999999

$ws.Range("H138").Value = 2655.8044
$ws.Range("J138").Value = 3181.4
$ws.Range("L138").Value = 9544.200000000001
$ws.Range("N138").Value = -19824.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5411.763
$ws.Range("I32").Value = 4517.6855
$ws.Range("K32").Value = 4517.6855
$ws.Range("M32").Value = -4230.6855

$ws.Range("H61").Value = 8994.272000000001
$ws.Range("I61").Value = 6142.125
$ws.Range("K61").Value = 6142.125
$ws.Range("M61").Value = -5930.125

$ws.Range("H88").Value = 3634.4285
$ws.Range("I88").Value = 3583.3333
$ws.Range("K88").Value = 3583.3333
$ws.Range("M88").Value = -3177.3333

$ws.Range("H91").Value = 3634.4285
$ws.Range("I91").Value = 3583.3333
$ws.Range("K91").Value = 3583.3333
$ws.Range("M91").Value = -2179.3333

$ws.Range("H102").Value = 2093.3333
$ws.Range("I102").Value = 2193.7273
$ws.Range("J102").Value = 989
$ws.Range("K102").Value = 2193.7273
$ws.Range("L102").Value = 989
$ws.Range("M102").Value = -571.7273
$ws.Range("N102").Value = -4233

$ws.Range("H122").Value = 2637.125
$ws.Range("I122").Value = 1925.8889
$ws.Range("K122").Value = 5777.6667
$ws.Range("M122").Value = -3327.6667

$ws.Range("H136").Value = 8994.272000000001
$ws.Range("I136").Value = 6142.125
$ws.Range("K136").Value = 18426.375
$ws.Range("M136").Value = -15876.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2932.5334
$ws.Range("I20").Value = 2083.6365
$ws.Range("J20").Value = 5267
$ws.Range("K20").Value = 2083.6365
$ws.Range("L20").Value = 5267
$ws.Range("M20").Value = -1836.6365
$ws.Range("N20").Value = -5761

$ws.Range("H59").Value = 89993
$ws.Range("J59").Value = 89993
$ws.Range("L59").Value = 89993
$ws.Range("N59").Value = -91687

$ws.Range("H86").Value = 4900.1763
$ws.Range("J86").Value = 9681
$ws.Range("L86").Value = 9681
$ws.Range("N86").Value = -11927

$ws.Range("H89").Value = 4900.1763
$ws.Range("J89").Value = 9681
$ws.Range("L89").Value = 48405
$ws.Range("N89").Value = -59637

$ws.Range("H94").Value = 1248.3529
$ws.Range("I94").Value = 1223
$ws.Range("K94").Value = 1223
$ws.Range("M94").Value = -772

$ws.Range("H105").Value = 8838.4375
$ws.Range("I105").Value = 3935.111
$ws.Range("K105").Value = 3935.111
$ws.Range("M105").Value = -2188.111

$ws.Range("H107").Value = 1042
$ws.Range("I107").Value = 1061.2354
$ws.Range("K107").Value = 1061.2354
$ws.Range("M107").Value = 858.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 8075
$ws.Range("J45").Value = 7600
$ws.Range("L45").Value = 7600
$ws.Range("N45").Value = -8786

$ws.Range("H48").Value = 27025.5
$ws.Range("J48").Value = 27025.5
$ws.Range("L48").Value = 27025.5
$ws.Range("N48").Value = -27977.5

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H122").Value = 4899.7666
$ws.Range("I122").Value = 1674
$ws.Range("J122").Value = 8586.357
$ws.Range("K122").Value = 5022
$ws.Range("L122").Value = 25759.071
$ws.Range("M122").Value = -2572
$ws.Range("N122").Value = -30659.071

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1800
$ws.Range("J51").Value = 1800
$ws.Range("L51").Value = 5400
$ws.Range("N51").Value = -6320

$ws.Range("H55").Value = 1998.8
$ws.Range("I55").Value = 996.3333
$ws.Range("J55").Value = 3502.5
$ws.Range("K55").Value = 2988.9999
$ws.Range("L55").Value = 10507.5
$ws.Range("M55").Value = -2811.9999
$ws.Range("N55").Value = -10861.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13759.808
$ws.Range("I70").Value = 5438
$ws.Range("K70").Value = 5438
$ws.Range("M70").Value = -5168

$ws.Range("H73").Value = 13759.808
$ws.Range("I73").Value = 5438
$ws.Range("K73").Value = 5438
$ws.Range("M73").Value = -4502

$ws.Range("H97").Value = 937.88464
$ws.Range("J97").Value = 1392.1
$ws.Range("L97").Value = 1392.1
$ws.Range("N97").Value = -2384.1

$ws.Range("H122").Value = 2899.4783
$ws.Range("I122").Value = 2594.375
$ws.Range("J122").Value = 3596.8572
$ws.Range("K122").Value = 7783.125
$ws.Range("L122").Value = 10790.5716
$ws.Range("M122").Value = -5333.125
$ws.Range("N122").Value = -15690.5716

$ws.Range("H132").Value = 57354.58
$ws.Range("I132").Value = 73979.21000000001
$ws.Range("J132").Value = 10805.6
$ws.Range("K132").Value = 221937.63
$ws.Range("L132").Value = 32416.8
$ws.Range("M132").Value = -219407.63
$ws.Range("N132").Value = -37476.8

$ws.Range("H136").Value = 38101.273
$ws.Range("J136").Value = 38101.273
$ws.Range("L136").Value = 114303.819
$ws.Range("N136").Value = -119403.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9004.071
$ws.Range("I7").Value = 3325.5557
$ws.Range("J7").Value = 19225.4
$ws.Range("K7").Value = 3325.5557
$ws.Range("L7").Value = 19225.4
$ws.Range("M7").Value = -3213.5557
$ws.Range("N7").Value = -19449.4

$ws.Range("H46").Value = 3045.875
$ws.Range("J46").Value = 3306.4614
$ws.Range("L46").Value = 3306.4614
$ws.Range("N46").Value = -3682.4614

$ws.Range("H55").Value = 1925088.6
$ws.Range("I55").Value = 3572290.8
$ws.Range("J55").Value = 3352.75
$ws.Range("K55").Value = 3572290.8
$ws.Range("L55").Value = 3352.75
$ws.Range("M55").Value = -3572117.8
$ws.Range("N55").Value = -3698.75

$ws.Range("H61").Value = 4629.5835
$ws.Range("I61").Value = 1194.25
$ws.Range("K61").Value = 1194.25
$ws.Range("M61").Value = -992.25

$ws.Range("H82").Value = 3534.7856
$ws.Range("I82").Value = 1454.5555
$ws.Range("K82").Value = 1454.5555
$ws.Range("M82").Value = -1093.5555

$ws.Range("H85").Value = 3534.7856
$ws.Range("I85").Value = 1454.5555
$ws.Range("K85").Value = 1454.5555
$ws.Range("M85").Value = -206.5554999999999

$ws.Range("H113").Value = 4629.5835
$ws.Range("I113").Value = 1194.25
$ws.Range("K113").Value = 1194.25
$ws.Range("M113").Value = 975.75

$ws.Range("H126").Value = 9004.071
$ws.Range("I126").Value = 3325.5557
$ws.Range("J126").Value = 19225.4
$ws.Range("K126").Value = 9976.667099999999
$ws.Range("L126").Value = 57676.2
$ws.Range("M126").Value = -7506.667099999999
$ws.Range("N126").Value = -62616.2

$ws.Range("H132").Value = 3572.3333
$ws.Range("I132").Value = 1340.125
$ws.Range("J132").Value = 8036.75
$ws.Range("K132").Value = 4020.375
$ws.Range("L132").Value = 24110.25
$ws.Range("M132").Value = -1490.375
$ws.Range("N132").Value = -29170.25

$ws.Range("H136").Value = 7540.909
$ws.Range("I136").Value = 4775.5
$ws.Range("J136").Value = 8577.9375
$ws.Range("K136").Value = 14326.5
$ws.Range("L136").Value = 25733.8125
$ws.Range("M136").Value = -11776.5
$ws.Range("N136").Value = -30833.8125

$ws.Range("H140").Value = 74268.60000000001
$ws.Range("J140").Value = 74085.75
$ws.Range("L140").Value = 74085.75
$ws.Range("N140").Value = -84445.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4876
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -12746

$ws.Range("H122").Value = 6391.9116
$ws.Range("I122").Value = 1741.2632
$ws.Range("J122").Value = 12282.733
$ws.Range("K122").Value = 5223.7896
$ws.Range("L122").Value = 36848.199
$ws.Range("M122").Value = -2773.7896
$ws.Range("N122").Value = -41748.199

$ws.Range("H126").Value = 3322.7
$ws.Range("I126").Value = 1897
$ws.Range("J126").Value = 5461.25
$ws.Range("K126").Value = 5691
$ws.Range("L126").Value = 16383.75
$ws.Range("M126").Value = -3221
$ws.Range("N126").Value = -21323.75

$ws.Range("H132").Value = 4969.1177
$ws.Range("I132").Value = 4552.1787
$ws.Range("K132").Value = 13656.5361
$ws.Range("M132").Value = -11126.5361
